$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, avoiding Excel auto-number/date coercion,
# while leaving the cell style/format untouched (matches source which has no "s" attr changes).
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "30.711.42"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "1.889.13"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue ($ws.Range("D5")) "247.89"
$ws.Range("E5").Value = "  +2.24%  "
Set-TextValue ($ws.Range("D6")) "1.001"
$ws.Range("E6").Value = "  +0.15%  "
Set-TextValue ($ws.Range("D7")) "0.4941"
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue ($ws.Range("D8")) "0.2962"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "1.888.39"
$ws.Range("E10").Value = "  +0.55%  "
Set-TextValue ($ws.Range("D11")) "17.20"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue ($ws.Range("D12")) "0.07238"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue ($ws.Range("D13")) "91.88"
$ws.Range("E13").Value = "  +6.38%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue ($ws.Range("D14")) "5.076"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue ($ws.Range("D15")) "0.6780"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "30.682.98"
$ws.Range("E16").Value = "  +2.37%  "
Set-TextValue ($ws.Range("D17")) "0.000007976"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("E18").Value = "  +0.08%  "
Set-TextValue ($ws.Range("D19")) "13.23"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "2.136.72"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").Value = "  +0.18%  "
Set-TextValue ($ws.Range("D22")) "4.835"
$ws.Range("E22").Value = "  +1.46%  "
Set-TextValue ($ws.Range("D23")) "189.42"
$ws.Range("E23").Value = "  +34.00%  "
Set-TextValue ($ws.Range("D24")) "6.063"
$ws.Range("E24").Value = "  +5.69%  "
Set-TextValue ($ws.Range("D25")) "9.371"
$ws.Range("E25").Value = "  +3.52%  "
Set-TextValue ($ws.Range("D26")) "156.62"
$ws.Range("E26").Value = "  +4.69%  "
Set-TextValue ($ws.Range("D27")) "19.05"
$ws.Range("E27").Value = "  +12.02%  "
Set-TextValue ($ws.Range("D28")) "1.910"
$ws.Range("E28").Value = "  -0.37%  "
Set-TextValue ($ws.Range("D29")) "1.404"
$ws.Range("E29").Value = "  +0.77%  "
Set-TextValue ($ws.Range("D30")) "4.310"
$ws.Range("E30").Value = "  +3.05%  "
Set-TextValue ($ws.Range("D31")) "0.09008"
$ws.Range("E31").Value = "  +3.20%  "
Set-TextValue ($ws.Range("D32")) "4.015"
$ws.Range("E32").Value = "  +1.43%  "
Set-TextValue ($ws.Range("D33")) "0.05190"
$ws.Range("E33").Value = "  +2.92%  "
Set-TextValue ($ws.Range("D34")) "0.7445"
$ws.Range("E34").Value = "  +4.83%  "
Set-TextValue ($ws.Range("D35")) "1.118"
$ws.Range("E35").Value = "  +0.20%  "
Set-TextValue ($ws.Range("D36")) "2.724"
$ws.Range("E36").Value = "  +2.13%  "
Set-TextValue ($ws.Range("D37")) "0.01842"
$ws.Range("E37").Value = "  +3.26%  "
Set-TextValue ($ws.Range("D38")) "2.674"
$ws.Range("E38").Value = "  -0.58%  "
Set-TextValue ($ws.Range("D39")) "2.162"
$ws.Range("E39").Value = "  -0.07%  "
Set-TextValue ($ws.Range("D40")) "0.9414"
$ws.Range("E40").Value = "  +1.12%  "
Set-TextValue ($ws.Range("D41")) "0.4425"
$ws.Range("E41").Value = "  +4.31%  "
Set-TextValue ($ws.Range("D42")) "105.79"
$ws.Range("E42").Value = "  +3.13%  "
Set-TextValue ($ws.Range("D43")) "1.001"
$ws.Range("E43").Value = "  +0.26%  "
Set-TextValue ($ws.Range("D44")) "5.765"
$ws.Range("E44").Value = "  +0.23%  "
Set-TextValue ($ws.Range("D45")) "7.633"
$ws.Range("E45").Value = "  +2.86%  "
Set-TextValue ($ws.Range("D46")) "0.1340"
$ws.Range("E46").Value = "  +5.69%  "
Set-TextValue ($ws.Range("D47")) "0.05842"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue ($ws.Range("D48")) "8.705"
$ws.Range("E48").Value = "  +4.70%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue ($ws.Range("D49")) "1.426"
$ws.Range("E49").Value = "  +7.04%  "
Set-TextValue ($ws.Range("D50")) "0.3943"
$ws.Range("E50").Value = "  +4.44%  "
Set-TextValue ($ws.Range("D51")) "33.55"
$ws.Range("E51").Value = "  +3.43%  "
